$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Revision list"
$ws2 = $wb.Worksheets.Item(2)   # "Definitions"

# --- Sheet2 "Definitions": fill in row 20 (ID 17 - antalPersoner) ---
$ws2.Range("B20").Value = "antalPersoner"
$ws2.Range("C20").Value = "antal af personer der bestiller den tur"
$ws2.Range("D20").Value = "all"
$ws2.Range("G20").Value = "antal passegerer : same but change to antalPersoner"

# Row 20 grows tall because of wrapped text across several lines (matches
# the heights already used on the other rows of this table).
$ws2.Rows.Item(20).RowHeight = 86.4

# --- Sheet1 "Revision list": append a new revision-history row ---
$ws1.Range("A14").Value = "Elaboration draft1"
$ws1.Range("B14").Value = "13.maj 2016"
$ws1.Range("C14").Value = "17 : antal personer "
$ws1.Range("D14").Value = "Jonas og Juyoung Choi"

# Column D on the revision list widens to fit the new description text.
$ws1.Columns.Item(4).ColumnWidth = 18.5

# --- Selections: restore per-sheet active cell, sheet1 stays the active tab ---
# (select on sheet2 first, then sheet1 last, so sheet1 ends up the active tab,
#  matching tabSelected="1" staying on "Revision list")
$ws2.Range("I19").Select() | Out-Null
$ws1.Range("F11").Select() | Out-Null
